$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the source diff.
# Values are written with a leading apostrophe to force text interpretation
# (preventing Excel from auto-converting numeric-looking strings into numbers,
# which would lose formatting such as trailing zeros or multi-dot grouping),
# then ClearFormats() removes the quote-prefix / text-format styling so the
# cell keeps the workbook's original (default) style.

$c = $ws.Range("D2")
$c.Value = "'28.962.31"
$c.ClearFormats()

$c = $ws.Range("E2")
$c.Value = "'  -2.25%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.Value = "'1.902.24"
$c.ClearFormats()

$c = $ws.Range("E3")
$c.Value = "'  -4.51%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.Value = "'  +0.40%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.Value = "'324.37"
$c.ClearFormats()

$c = $ws.Range("E5")
$c.Value = "'  -1.36%  "
$c.ClearFormats()

$c = $ws.Range("E6")
$c.Value = "'  +0.37%  "
$c.ClearFormats()

$c = $ws.Range("D7")
$c.Value = "'0.4593"
$c.ClearFormats()

$c = $ws.Range("E7")
$c.Value = "'  -1.81%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.Value = "'0.3815"
$c.ClearFormats()

$c = $ws.Range("E8")
$c.Value = "'  -3.52%  "
$c.ClearFormats()

$c = $ws.Range("E9")
$c.Value = "'  -2.07%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.Value = "'0.07701"
$c.ClearFormats()

$c = $ws.Range("E10")
$c.Value = "'  -4.67%  "
$c.ClearFormats()

$c = $ws.Range("D11")
$c.Value = "'0.9801"
$c.ClearFormats()

$c = $ws.Range("E11")
$c.Value = "'  -2.10%  "
$c.ClearFormats()

$c = $ws.Range("E12")
$c.Value = "'  -4.06%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.Value = "'1.922.23"
$c.ClearFormats()

$c = $ws.Range("E13")
$c.Value = "'  -3.23%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.Value = "'6.967"
$c.ClearFormats()

$c = $ws.Range("E14")
$c.Value = "'  -3.79%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.Value = "'5.668"
$c.ClearFormats()

$c = $ws.Range("E15")
$c.Value = "'  -3.43%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.Value = "'0.07059"
$c.ClearFormats()

$c = $ws.Range("E16")
$c.Value = "'  -1.00%  "
$c.ClearFormats()

$c = $ws.Range("D17")
$c.Value = "'1.005"
$c.ClearFormats()

$c = $ws.Range("E17")
$c.Value = "'  +0.32%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.Value = "'83.86"
$c.ClearFormats()

$c = $ws.Range("E18")
$c.Value = "'  -5.62%  "
$c.ClearFormats()

$c = $ws.Range("D19")
$c.Value = "'0.000009528"
$c.ClearFormats()

$c = $ws.Range("E19")
$c.Value = "'  -5.01%  "
$c.ClearFormats()

$c = $ws.Range("D20")
$c.Value = "'16.73"
$c.ClearFormats()

$c = $ws.Range("E20")
$c.Value = "'  -4.02%  "
$c.ClearFormats()

$c = $ws.Range("E21")
$c.Value = "'  +0.40%  "
$c.ClearFormats()

$c = $ws.Range("D22")
$c.Value = "'29.000.19"
$c.ClearFormats()

$c = $ws.Range("E22")
$c.Value = "'  -2.06%  "
$c.ClearFormats()

$c = $ws.Range("D23")
$c.Value = "'5.328"
$c.ClearFormats()

$c = $ws.Range("E23")
$c.Value = "'  -4.18%  "
$c.ClearFormats()

$c = $ws.Range("E24")
$c.Value = "'  -3.18%  "
$c.ClearFormats()

$c = $ws.Range("D25")
$c.Value = "'2.205.04"
$c.ClearFormats()

$c = $ws.Range("E25")
$c.Value = "'  -1.08%  "
$c.ClearFormats()

$c = $ws.Range("D26")
$c.Value = "'2.101"
$c.ClearFormats()

$c = $ws.Range("E26")
$c.Value = "'  -0.38%  "
$c.ClearFormats()

$c = $ws.Range("D27")
$c.Value = "'157.37"
$c.ClearFormats()

$c = $ws.Range("E27")
$c.Value = "'  -0.25%  "
$c.ClearFormats()

$c = $ws.Range("E28")
$c.Value = "'  -2.88%  "
$c.ClearFormats()

$c = $ws.Range("D29")
$c.Value = "'5.591"
$c.ClearFormats()

$c = $ws.Range("E29")
$c.Value = "'  -6.57%  "
$c.ClearFormats()

$c = $ws.Range("D30")
$c.Value = "'117.56"
$c.ClearFormats()

$c = $ws.Range("E30")
$c.Value = "'  -2.39%  "
$c.ClearFormats()

$c = $ws.Range("D31")
$c.Value = "'1.851"
$c.ClearFormats()

$c = $ws.Range("E31")
$c.Value = "'  -4.82%  "
$c.ClearFormats()

$c = $ws.Range("D32")
$c.Value = "'0.09283"
$c.ClearFormats()

$c = $ws.Range("E32")
$c.Value = "'  -1.73%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.Value = "'0.8600"
$c.ClearFormats()

$c = $ws.Range("E33")
$c.Value = "'  -6.11%  "
$c.ClearFormats()

$c = $ws.Range("D34")
$c.Value = "'5.079"
$c.ClearFormats()

$c = $ws.Range("E34")
$c.Value = "'  -3.87%  "
$c.ClearFormats()

$c = $ws.Range("D35")
$c.Value = "'1.249"
$c.ClearFormats()

$c = $ws.Range("E35")
$c.Value = "'  -7.83%  "
$c.ClearFormats()

$c = $ws.Range("D36")
$c.Value = "'3.014"
$c.ClearFormats()

$c = $ws.Range("E36")
$c.Value = "'  -5.35%  "
$c.ClearFormats()

$c = $ws.Range("D37")
$c.Value = "'0.05687"
$c.ClearFormats()

$c = $ws.Range("E37")
$c.Value = "'  -2.91%  "
$c.ClearFormats()

$c = $ws.Range("D38")
$c.Value = "'1.147"
$c.ClearFormats()

$c = $ws.Range("E38")
$c.Value = "'  -2.67%  "
$c.ClearFormats()

$c = $ws.Range("D39")
$c.Value = "'1.003"
$c.ClearFormats()

$c = $ws.Range("E39")
$c.Value = "'  +0.27%  "
$c.ClearFormats()

$c = $ws.Range("D40")
$c.Value = "'0.02035"
$c.ClearFormats()

$c = $ws.Range("E40")
$c.Value = "'  -4.41%  "
$c.ClearFormats()

$c = $ws.Range("D41")
$c.Value = "'7.464"
$c.ClearFormats()

$c = $ws.Range("E41")
$c.Value = "'  -5.63%  "
$c.ClearFormats()

$c = $ws.Range("D42")
$c.Value = "'0.5513"
$c.ClearFormats()

$c = $ws.Range("E42")
$c.Value = "'  -4.74%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.Value = "'0.1754"
$c.ClearFormats()

$c = $ws.Range("E43")
$c.Value = "'  -3.95%  "
$c.ClearFormats()

$c = $ws.Range("D44")
$c.Value = "'9.265"
$c.ClearFormats()

$c = $ws.Range("E44")
$c.Value = "'  -6.45%  "
$c.ClearFormats()

$c = $ws.Range("D45")
$c.Value = "'2.748"
$c.ClearFormats()

$c = $ws.Range("E45")
$c.Value = "'  -2.59%  "
$c.ClearFormats()

$c = $ws.Range("D46")
$c.Value = "'0.5191"
$c.ClearFormats()

$c = $ws.Range("D47")
$c.Value = "'11.26"
$c.ClearFormats()

$c = $ws.Range("E47")
$c.Value = "'  -7.06%  "
$c.ClearFormats()

$c = $ws.Range("D48")
$c.Value = "'2.084"
$c.ClearFormats()

$c = $ws.Range("E48")
$c.Value = "'  -5.79%  "
$c.ClearFormats()

$c = $ws.Range("D49")
$c.Value = "'0.06821"
$c.ClearFormats()

$c = $ws.Range("E49")
$c.Value = "'  -2.23%  "
$c.ClearFormats()

$c = $ws.Range("D50")
$c.Value = "'111.28"
$c.ClearFormats()

$c = $ws.Range("E50")
$c.Value = "'  -2.48%  "
$c.ClearFormats()

$c = $ws.Range("D51")
$c.Value = "'1.777"
$c.ClearFormats()

$c = $ws.Range("E51")
$c.Value = "'  -5.10%  "
$c.ClearFormats()
